# edit.ps1 - apply 'Added v1.2 of valve block' changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate for the data area (keeps column widths/formats).
$ws.Cells.ClearContents()

# --- Pre-seed the shared-string table in the exact order the target file uses ---
# (Excel's sharedStrings.xml order follows write order; scratch cells are cleared
# afterwards but the strings remain interned because the real cells below reuse them.)
$ws.Range("Z1").Value = "depth neg"
$ws.Range("Z2").Value = "top corridor"
$ws.Range("Z3").Value = "bottom corridor"
$ws.Range("Z4").Value = "7 holes"
$ws.Range("Z5").Value = "15 holes"
$ws.Range("Z6").Value = "no holes"
$ws.Range("Z7").Value = "2 x 6 slots"
$ws.Range("Z8").Value = "2 x 13 holes"
$ws.Range("Z9").Value = "number of channels in bottom"
$ws.Range("Z10").Value = "individual channel width in bottom"
$ws.Range("Z11").Value = "area inside hosetail 16 mm"
$ws.Range("Z12").Value = "pressure controlled valve type a"
$ws.Range("Z13").Value = "pressure controlled valve type b"
$ws.Range("Z14").Value = "pressure controlled valve type c"
$ws.Range("Z15").Value = "one way valve type b"
$ws.Range("Z16").Value = "one way valve type c"
$ws.Range("Z17").Value = "area inside hosetail 19.5 mm"
$ws.Range("Z18").Value = "area inside hosetail 13.5 mm"
$ws.Range("Z19").Value = "mini one way valve"
$ws.Range("Z20").Value = "<-"
$ws.Range("Z21").Value = "oval via"

# --- Column widths (A widened, new C column added) ---
$ws.Columns("A").ColumnWidth = 32.666666666666664
$ws.Columns("C").ColumnWidth = 12

# --- Row content ---
$ws.Range("A3").Value = "depth neg"
$ws.Range("B3").Value = 12

$ws.Range("A4").Value = "number of channels in bottom"
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = "individual channel width in bottom"
$ws.Range("B5").Value = 9

$ws.Range("A7").Value = "area inside hosetail 19.5 mm"
$ws.Range("B7").Formula = "=PI()/4*(16.5*16.5)"

$ws.Range("A8").Value = "area inside hosetail 16 mm"
$ws.Range("B8").Formula = "=PI()/4*(13*13)"

$ws.Range("A9").Value = "area inside hosetail 13.5 mm"
$ws.Range("B9").Formula = "=PI()/4*(10.5*10.5)"

$ws.Range("A10").Value = "pressure controlled valve type a"
$ws.Range("B10").Formula = "=PI()/4*12*12+12*12"
$ws.Range("C10").Value = "no holes"

$ws.Range("A11").Value = "pressure controlled valve type b"
$ws.Range("B11").Formula = "=13*(PI()/4*3*3)"
$ws.Range("C11").Value = "2 x 13 holes"

$ws.Range("A12").Value = "pressure controlled valve type c"
$ws.Range("B12").Formula = "=0.6*B10"
$ws.Range("C12").Value = "2 x 6 slots"
$ws.Range("D12").Value = "<-"

$ws.Range("A13").Value = "one way valve type b"
$ws.Range("B13").Formula = "=15*(PI()/4*3*3)"
$ws.Range("C13").Value = "15 holes"

$ws.Range("A14").Value = "one way valve type c"
$ws.Range("B14").Formula = "=6*PI()/4*(5*5) + PI()/4*(4.5*4.5)"
$ws.Range("C14").Value = "7 holes"
$ws.Range("D14").Value = "<-"

$ws.Range("A15").Value = "oval via"
$ws.Range("B15").Formula = "=PI()/4*10*10+10*7"

$ws.Range("A16").Value = "mini one way valve"
$ws.Range("B16").Formula = "=2*(PI()/4*4.5*4.5)"

$ws.Range("A17").Value = "top corridor"
$ws.Range("B17").Formula = "=0.5*(17*(15-6.5))+(17*6.5)"

$ws.Range("A18").Value = "bottom corridor"
$ws.Range("B18").Formula = "=`$B`$4*(`$B`$5*(`$B`$3-`$B`$5/2)+0.125*`$B`$5*`$B`$5)-2*0.5*0.5*0.5"

# --- Clear the scratch cells now that the real cells above reuse those strings ---
$ws.Range("Z1:Z21").ClearContents()

# --- Restore the view's active cell/selection to match the authored state ---
[void]$ws.Range("A20").Select()

Write-Output "edit applied"
